$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 4).Value = '28.343.41'
$ws.Cells.Item(2, 5).Value = '  -0.86%  '
$ws.Cells.Item(3, 4).Value = '1.553.02'
$ws.Cells.Item(3, 5).Value = '  -1.49%  '
$ws.Cells.Item(4, 5).Value = '  -0.30%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '210.24'
$ws.Cells.Item(5, 5).Value = '  -1.42%  '
$ws.Cells.Item(6, 5).Value = '  -2.12%  '
$ws.Cells.Item(7, 5).Value = '  -0.27%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '23.86'
$ws.Cells.Item(8, 5).Value = '  -0.53%  '
$ws.Cells.Item(9, 5).Value = '  -2.07%  '
$ws.Cells.Item(10, 5).Value = '  -1.41%  '
$ws.Cells.Item(11, 5).Value = '  -0.15%  '
$ws.Cells.Item(12, 4).Value = '1.774.79'
$ws.Cells.Item(12, 5).Value = '  -1.55%  '
$ws.Cells.Item(13, 4).Value = '1.540.19'
$ws.Cells.Item(13, 5).Value = '  -2.15%  '
$ws.Cells.Item(14, 4).Value = '28.334.67'
$ws.Cells.Item(14, 5).Value = '  -0.92%  '
$ws.Cells.Item(15, 5).Value = '  -1.91%  '
$ws.Cells.Item(16, 5).Value = '  -1.79%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '60.89'
$ws.Cells.Item(17, 5).Value = '  -2.20%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '228.04'
$ws.Cells.Item(18, 5).Value = '  -1.10%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '7.36'
$ws.Cells.Item(19, 5).Value = '  -0.35%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0675'
$ws.Cells.Item(20, 5).Value = '  -2.26%  '
$ws.Cells.Item(21, 5).Value = '  -0.16%  '
$ws.Cells.Item(22, 5).Value = '  +1.01%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '8.94'
$ws.Cells.Item(23, 5).Value = '  -2.42%  '
$ws.Cells.Item(24, 5).Value = '  -2.48%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '151.06'
$ws.Cells.Item(25, 5).Value = '  -0.29%  '
$ws.Cells.Item(26, 5).Value = '  -1.69%  '
$ws.Cells.Item(27, 5).Value = '  -1.15%  '
$ws.Cells.Item(28, 5).Value = '  -0.20%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '6.25'
$ws.Cells.Item(29, 5).Value = '  -3.00%  '
$ws.Cells.Item(30, 5).Value = '  -3.00%  '
$ws.Cells.Item(31, 5).Value = '  -4.52%  '
$ws.Cells.Item(32, 5).Value = '  -1.21%  '
$ws.Cells.Item(33, 4).Value = '1.389.56'
$ws.Cells.Item(33, 5).Value = '  -0.75%  '
$ws.Cells.Item(34, 5).Value = '  -2.59%  '
$ws.Cells.Item(35, 5).Value = '  +2.37%  '
$ws.Cells.Item(36, 5).Value = '  -3.60%  '
$ws.Cells.Item(37, 5).Value = '  -0.98%  '
$ws.Cells.Item(38, 5).Value = '  -1.46%  '
$ws.Cells.Item(39, 5).Value = '  -2.75%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.513'
$ws.Cells.Item(40, 5).Value = '  -1.89%  '
$ws.Cells.Item(41, 5).Value = '  +1.40%  '
$ws.Cells.Item(42, 5).Value = '  -0.22%  '
$ws.Cells.Item(43, 5).Value = '  -2.06%  '
$ws.Cells.Item(44, 5).Value = '  -1.59%  '
$ws.Cells.Item(45, 5).Value = '  -1.95%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '61.95'
$ws.Cells.Item(46, 5).Value = '  -1.54%  '
$ws.Cells.Item(47, 4).Value = '1.687.53'
$ws.Cells.Item(47, 5).Value = '  -1.62%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.868'
$ws.Cells.Item(48, 5).Value = '  -9.73%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '85.63'
$ws.Cells.Item(49, 5).Value = '  -0.98%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '42.83'
$ws.Cells.Item(50, 5).Value = '  +7.90%  '
$ws.Cells.Item(51, 5).Value = '  +0.09%  '
